$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

# Update the "Lower Right Cell" values for the extra scenario blocks
# (row range extended from *43 to *46)
$ws.Range("D5").Value = "A46"
$ws.Range("D6").Value = "B46"
$ws.Range("D7").Value = "C46"
$ws.Range("D8").Value = "G46"
$ws.Range("D9").Value = "H46"
$ws.Range("D10").Value = "I46"
$ws.Range("D11").Value = "J46"

# Update the selection to match the new active range
$ws.Range("D5:D11").Select()
